$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Report the Monday hours (7) for the week commencing 43234 (row 18)
$ws.Range("B18").Value = 7

# Reflect the final active selection left by the user after entering data
$ws.Range("D20").Select()
